$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$ws.Range("H3").Value = -104
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "04-Nov-2025"
$ws.Range("I3").NumberFormat = "general"

$ws.Range("H4").Value = 699
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "04-Nov-2025"
$ws.Range("I4").NumberFormat = "general"
